$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($rangeAddr, $val) {
    $c = $ws.Range($rangeAddr)
    $s = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $s
}

Set-TextValue "A2" "2025-06-11"
$ws.Range("B2").Value = 2
Set-TextValue "C2" "BEMOL S/A"
Set-TextValue "D2" "357349"
$ws.Range("E2").Value = 13546
Set-TextValue "F2" "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON"
$ws.Range("G2").Value = -291
$ws.Range("H2").Value = 1.09
$ws.Range("I2").Value = 0.3

Set-TextValue "A3" "2025-06-11"
$ws.Range("B3").Value = 2
Set-TextValue "C3" "BEMOL S/A"
Set-TextValue "D3" "357392"
$ws.Range("E3").Value = 4408
Set-TextValue "F3" "RING LIGHT 10 POLEGADAS COM TRIPE"
$ws.Range("G3").Value = -450
$ws.Range("H3").Value = 1.01
$ws.Range("I3").Value = 0.15

Set-TextValue "A4" "2025-06-11"
$ws.Range("B4").Value = 2
Set-TextValue "C4" "BEMOL S/A"
Set-TextValue "D4" "357402"
$ws.Range("E4").Value = 13546
Set-TextValue "F4" "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON"
$ws.Range("G4").Value = -291
$ws.Range("H4").Value = 1.09
$ws.Range("I4").Value = 0.3

Set-TextValue "A5" "2025-06-12"
$ws.Range("B5").Value = 2
Set-TextValue "C5" "BEMOL S/A"
Set-TextValue "D5" "358537"
$ws.Range("E5").Value = 13588
Set-TextValue "F5" "CANETA STYLLUS ACTIVA AGOLD"
$ws.Range("G5").Value = -74
$ws.Range("H5").Value = 1.05
$ws.Range("I5").Value = 0.22

Set-TextValue "A6" "2025-06-12"
$ws.Range("B6").Value = 2
Set-TextValue "C6" "BEMOL S/A"
Set-TextValue "D6" "358540"
$ws.Range("E6").Value = 13546
Set-TextValue "F6" "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON"
$ws.Range("G6").Value = -291
$ws.Range("H6").Value = 1.09
$ws.Range("I6").Value = 0.3

Set-TextValue "A7" "2025-06-13"
$ws.Range("B7").Value = 2
Set-TextValue "C7" "BEMOL S/A"
Set-TextValue "D7" "359654"
$ws.Range("E7").Value = 12680
Set-TextValue "F7" "BASTAO DE LUZ RGB LED TOMATE"
$ws.Range("G7").Value = -28
$ws.Range("H7").Value = 1.08
$ws.Range("I7").Value = 0.27

Set-TextValue "A8" "2025-06-14"
$ws.Range("B8").Value = 2
Set-TextValue "C8" "SOCIEDADE MICHELIN DE PARTICIPACOES INDUST E COMERCIO LTDA"
Set-TextValue "D8" "360346"
$ws.Range("E8").Value = 13588
Set-TextValue "F8" "CANETA STYLLUS ACTIVA AGOLD"
$ws.Range("G8").Value = -74
$ws.Range("H8").Value = 1.05
$ws.Range("I8").Value = 0.22

Set-TextValue "A9" "2025-06-16"
$ws.Range("B9").Value = 2
Set-TextValue "C9" "BEMOL S/A"
Set-TextValue "D9" "362396"
$ws.Range("E9").Value = 13079
Set-TextValue "F9" "FONE BLUETOOTH BASIKE TWS FON6694"
$ws.Range("G9").Value = -434
$ws.Range("H9").Value = 1.08
$ws.Range("I9").Value = 0.29

Set-TextValue "A10" "2025-06-16"
$ws.Range("B10").Value = 2
Set-TextValue "C10" "BEMOL S/A"
Set-TextValue "D10" "362404"
$ws.Range("E10").Value = 13546
Set-TextValue "F10" "FONE SEM FIO BOX PRETO ESTEREO EARBUDS CASE CARREGADOR BLUETOOTH LETRON"
$ws.Range("G10").Value = -291
$ws.Range("H10").Value = 1.09
$ws.Range("I10").Value = 0.3

Set-TextValue "A11" "2025-06-19"
$ws.Range("B11").Value = 2
Set-TextValue "C11" "BEMOL S/A"
Set-TextValue "D11" "365782"
$ws.Range("E11").Value = 13079
Set-TextValue "F11" "FONE BLUETOOTH BASIKE TWS FON6694"
$ws.Range("G11").Value = -434
$ws.Range("H11").Value = 1.08
$ws.Range("I11").Value = 0.29

Set-TextValue "A12" "2025-06-19"
$ws.Range("B12").Value = 2
Set-TextValue "C12" "BEMOL S/A"
Set-TextValue "D12" "366707"
$ws.Range("E12").Value = 13079
Set-TextValue "F12" "FONE BLUETOOTH BASIKE TWS FON6694"
$ws.Range("G12").Value = -434
$ws.Range("H12").Value = 1.08
$ws.Range("I12").Value = 0.29

Set-TextValue "A13" "2025-06-24"
$ws.Range("B13").Value = 2
Set-TextValue "C13" "BEMOL S/A"
Set-TextValue "D13" "370495"
$ws.Range("E13").Value = 46217
Set-TextValue "F13" "SMART WATCH HMASTON INK12"
$ws.Range("G13").Value = -88
$ws.Range("H13").Value = 1.03
$ws.Range("I13").Value = 0.18
